# Generate Report for Handoff
#
# Replaces the two in-flight handoff records (formerly handed-back /
# in-sync with en-US) with the new "ready for handoff" batch: new file
# GUIDs, a combined target .xlf per locale, and refreshed timestamps.
# The now-unused "Latest Target File" / "Latest Handback File" columns
# (F/G) on the locale sheets are cleared since there is no handback yet.

$wb = $excel.ActiveWorkbook

# ---- new identifiers coming from the handoff batch -----------------
$guid1 = "6ad202f5-6023-44f2-8fe0-3e95703b1899"
$guid2 = "ffff74c1378f-cd0e-46f7-89ac-d20417659548"
$xlfHash = "57ba641c50aade9f316b9f411921c77b0b4ddfe1"

$status = "Ready for handoff"
$overviewDate = "2016-41-14 04:41:20"
$zhHandoffDatetime = "2016-03-14 04:41:17"
$deHandoffDatetime = "2016-03-14 04:41:20"
$blankTarget = "0001-01-01 00:00:00"

$zhXlf = "$guid1.$xlfHash.zh-cn.xlf"
$deXlf = "$guid1.$xlfHash.de-de.xlf"

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/4aea4ddd05fe8de8876f226dccb38ed12255bc1d/e2e"
$zhXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f59929daa2857e6300ba8bb9f1034f09754f9954/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deXlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/48d96f72d98687b05574148713e09ac661473f9d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$guid1MdUrl = "$mdUrlBase/$guid1.md"
$guid2MdUrl = "$mdUrlBase/$guid2.md"
$zhXlfUrl = "$zhXlfUrlBase/$zhXlf"
$deXlfUrl = "$deXlfUrlBase/$deXlf"

# =====================================================================
# Overview sheet
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "$guid1.md"
$wsOverview.Range("B2").Value = $status
$wsOverview.Range("C2").Value = $status
$wsOverview.Range("D2").Value = $overviewDate

$wsOverview.Range("A3").Value = "$guid2.md"
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = $overviewDate

$wsOverview.Range("A2").Hyperlinks.Add($wsOverview.Range("A2"), $guid1MdUrl, [Type]::Missing, [Type]::Missing, "$guid1.md")
$wsOverview.Range("A3").Hyperlinks.Add($wsOverview.Range("A3"), $guid2MdUrl, [Type]::Missing, [Type]::Missing, "$guid2.md")

# =====================================================================
# zh-cn sheet
# =====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Range("A2").Value = "$guid1.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $status
$wsZh.Range("D2").Value = $zhXlf
$wsZh.Range("E2").Value = $zhHandoffDatetime
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("H2").Value = $blankTarget
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = "$guid2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $zhXlf
$wsZh.Range("E3").Value = $zhHandoffDatetime
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()
$wsZh.Range("H3").Value = $blankTarget
$wsZh.Range("I3").Value = "Include"

$wsZh.Range("A2").Hyperlinks.Add($wsZh.Range("A2"), $guid1MdUrl, [Type]::Missing, [Type]::Missing, "$guid1.md")
$wsZh.Range("B2").Hyperlinks.Add($wsZh.Range("B2"), $guid1MdUrl, [Type]::Missing, [Type]::Missing, ".md")
$wsZh.Range("D2").Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlf)

$wsZh.Range("A3").Hyperlinks.Add($wsZh.Range("A3"), $guid2MdUrl, [Type]::Missing, [Type]::Missing, "$guid2.md")
$wsZh.Range("B3").Hyperlinks.Add($wsZh.Range("B3"), $guid2MdUrl, [Type]::Missing, [Type]::Missing, ".md")
$wsZh.Range("D3").Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlf)

# =====================================================================
# de-de sheet
# =====================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Range("A2").Value = "$guid1.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $status
$wsDe.Range("D2").Value = $deXlf
$wsDe.Range("E2").Value = $deHandoffDatetime
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("H2").Value = $blankTarget
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = "$guid2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $deXlf
$wsDe.Range("E3").Value = $deHandoffDatetime
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
$wsDe.Range("H3").Value = $blankTarget
$wsDe.Range("I3").Value = "Include"

$wsDe.Range("A2").Hyperlinks.Add($wsDe.Range("A2"), $guid1MdUrl, [Type]::Missing, [Type]::Missing, "$guid1.md")
$wsDe.Range("B2").Hyperlinks.Add($wsDe.Range("B2"), $guid1MdUrl, [Type]::Missing, [Type]::Missing, ".md")
$wsDe.Range("D2").Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlf)

$wsDe.Range("A3").Hyperlinks.Add($wsDe.Range("A3"), $guid2MdUrl, [Type]::Missing, [Type]::Missing, "$guid2.md")
$wsDe.Range("B3").Hyperlinks.Add($wsDe.Range("B3"), $guid2MdUrl, [Type]::Missing, [Type]::Missing, ".md")
$wsDe.Range("D3").Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlf)
